$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fix the typo in the existing "Can no iterate..." note (row 6, column E)
$ws.Range("E6").Value = "Can now iterate through the individual sheet and print to pdf... Understand how highlighting works"

# Row 7: add date, start time, end time, and note "Formatting the PDF"
$ws.Range("A6:C6").Copy()
$ws.Range("A7:C7").PasteSpecial(-4122)
$ws.Range("A7").Value = 41815
$ws.Range("B7").Value = 0.83333333333333337
$ws.Range("C7").Value = 0.875
$ws.Range("E7").Value = "Formatting the PDF"

# Row 8: add date, start time, end time
$ws.Range("A6:C6").Copy()
$ws.Range("A8:C8").PasteSpecial(-4122)
$ws.Range("A8").Value = 41816
$ws.Range("B8").Value = 0.80555555555555547
$ws.Range("C8").Value = 0.84722222222222221

# Update the active selection cell shown in the saved view
$ws.Range("E8").Select()
